$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New summary cell under column J (|S*|/n) -> average of the 10 data rows
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New summary rows 14-17 with labels in column A and aggregate formulas in column B
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Bold + slightly larger, vertically centered font for the new summary values
foreach ($addr in @("B14", "B15", "B16", "B17")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.Font.Size = 12
    $cell.VerticalAlignment = -4108
}

# Bold font for the column-J average cell
$ws.Range("J12").Font.Bold = $true

# Selection left on J12 (matches the cursor position saved with the workbook)
$ws.Range("J12").Select()

# Page setup tweaks made when the workbook was re-saved (A4, portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
